# The workbook gained one new data row: a new record for Apio
# (Macroferia Regional de Talca) was inserted right above the row that
# used to be row 83, pushing that row and everything below it down by
# one. The new row repeats the same market/category/quality/unit data as
# the (old) row 83, but carries a new date (serial 44757 = 2022-07-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 83 (and everything below it) down by one row.
$ws.Rows("83").Insert()

# Populate the newly-opened row 83 with the new record.
$ws.Cells.Item(83, 1).Value = 5
$ws.Cells.Item(83, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(83, 3).Value = "Maule"
$ws.Cells.Item(83, 4).Value = 44757
$ws.Cells.Item(83, 5).Value = 7
$ws.Cells.Item(83, 6).Value = 100112017
$ws.Cells.Item(83, 7).Value = "Apio"
$ws.Cells.Item(83, 8).Value = "Americana (o)"
$ws.Cells.Item(83, 9).Value = "Primera"
$ws.Cells.Item(83, 10).Value = 500
$ws.Cells.Item(83, 11).Value = 9000
$ws.Cells.Item(83, 12).Value = 9000
$ws.Cells.Item(83, 13).Value = 9000
$ws.Cells.Item(83, 14).Value = '$/docena de matas'
$ws.Cells.Item(83, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(83, 16).Value = 1500
$ws.Cells.Item(83, 17).Value = 6
$ws.Cells.Item(83, 18).Value = "Hortaliza"
